$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts existing rows 13-44 down to 14-45)
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the latest weekly price record
$ws.Cells.Item(13, 1).Value = 6
$ws.Cells.Item(13, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44804
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = 100112035
$ws.Cells.Item(13, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 310
$ws.Cells.Item(13, 11).Value = 17000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 17581
$ws.Cells.Item(13, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 16).Value = 1172
$ws.Cells.Item(13, 17).Value = 15
$ws.Cells.Item(13, 18).Value = "Hortaliza"
